# Sync automático del tracker (cada 3h)
# Appends 9 new match rows (480-488) to Sheet1, matching the incoming
# feed entries. event_id (col A) and fecha (col B) are fed as text,
# cuota (col F) stays numeric; resultado/profit (G/H) are left blank
# until the match is settled, same as the most recent existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Add-TrackerRow($Row, $EventId, $Fecha, $JugadorA, $JugadorB, $Pronostico, $Cuota) {
    # Force text storage for the id/date columns so they don't get
    # silently coerced into numbers / Excel serial dates.
    $ws.Cells.Item($Row, 1).NumberFormat = "@"
    $ws.Cells.Item($Row, 1).Value = $EventId

    $ws.Cells.Item($Row, 2).NumberFormat = "@"
    $ws.Cells.Item($Row, 2).Value = $Fecha

    $ws.Cells.Item($Row, 3).Value = $JugadorA
    $ws.Cells.Item($Row, 4).Value = $JugadorB
    $ws.Cells.Item($Row, 5).Value = $Pronostico

    $ws.Cells.Item($Row, 6).Value = $Cuota

    # resultado / profit: still undecided for these freshly-synced
    # matches, so they stay empty (same as the newest rows already in
    # the tracker).
    $ws.Cells.Item($Row, 7).Value = ""
    $ws.Cells.Item($Row, 8).Value = ""
}

Add-TrackerRow 480 "14494924" "2025-08-30" "Lorenzo Musetti"      "Flavio Cobolli"    "Gana Flavio Cobolli"      3.2
Add-TrackerRow 481 "14494922" "2025-08-30" "Jannik Sinner"        "Denis Shapovalov"  "Gana Denis Shapovalov"    21
Add-TrackerRow 482 "14494923" "2025-08-30" "Jaume Munar"          "Zizou Bergs"       "Gana Zizou Bergs"         2.75
Add-TrackerRow 483 "14495029" "2025-08-30" "Ekaterina Alexandrova" "Laura Siegemund"  "Gana Laura Siegemund"     4
Add-TrackerRow 484 "14495032" "2025-08-31" "Beatriz Haddad Maia"  "Maria Sakkari"     "Gana Beatriz Haddad Maia" 2.3
Add-TrackerRow 485 "14494924" "2025-08-30" "Lorenzo Musetti"      "Flavio Cobolli"    "Gana Flavio Cobolli"      3.2
Add-TrackerRow 486 "14494922" "2025-08-30" "Jannik Sinner"        "Denis Shapovalov"  "Gana Denis Shapovalov"    21
Add-TrackerRow 487 "14495029" "2025-08-30" "Ekaterina Alexandrova" "Laura Siegemund"  "Gana Laura Siegemund"     4
Add-TrackerRow 488 "14487604" "2025-08-30" "Daniel Rincon"        "Harold Mayot"      "Gana Daniel Rincon"       2.38
